# "some pictures and tables were updated"
# - add new worksheet "List1" with a summary table (Table 1)
# - tidy up dev.char sheet (extra column, selection) and mean.dev.time selection
# - scroll the workbook so the first visible tab is "dev.char"

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Add the new "List1" worksheet at the end of the tab strip
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws.Name = "List1"

# column widths
$ws.Columns.Item(1).ColumnWidth = 8.28515625
$ws.Columns.Item(2).ColumnWidth = 18.28515625
$ws.Columns.Item(3).ColumnWidth = 7.85546875
$ws.Columns.Item(4).ColumnWidth = 4
$ws.Columns.Item(5).ColumnWidth = 10.140625
$ws.Columns.Item(6).ColumnWidth = 16.7109375
$ws.Columns.Item(7).ColumnWidth = 12.5703125

# ---- data entry (order matters for shared-string table layout) ----
$ws.Range("B3").Value = "Temperature range"
$ws.Range("C3").Value = "R2"
$ws.Range("D3").Value = "Df"
$ws.Range("E3").Value = "p value"
$ws.Range("A4").Value = "egg"
$ws.Range("F3").Value = "k "
$ws.Range("B8").Value = "15-21"
$ws.Range("B4").Value = "15-25"
$ws.Range("A2").Value = "Table 1: Summary of development constants for S. watsoni for five developmental stages  (sum of effective temperatures (k) and lower developmental threshold (t))  (means and standard errors)."

$ws.Range("A3").Value = "Stage"
$ws.Range("G3").Value = "t"

$ws.Range("A5").Value = "L1"
$ws.Range("A6").Value = "L2"
$ws.Range("A7").Value = "L3"
$ws.Range("A8").Value = "Pupae"

$ws.Range("B5").Value = "15-25"
$ws.Range("B6").Value = "15-25"
$ws.Range("B7").Value = "15-25"

$ws.Range("C4").Value = 0.8134
$ws.Range("C5").Value = 0.9375
$ws.Range("C6").Value = 0.8768
$ws.Range("C7").Value = 0.8199
$ws.Range("C8").Value = 0.8563

$ws.Range("D4").Value = 220
$ws.Range("D5").Value = 171
$ws.Range("D6").Value = 206
$ws.Range("D7").Value = 27
$ws.Range("D8").Value = 10

$ws.Range("E4").Value = [double]"2.2e-16"
$ws.Range("E5").Value = [double]"2.2e-16"
$ws.Range("E6").Value = [double]"2.2e-16"
$ws.Range("E7").Value = [double]"1.486e-11"
$ws.Range("E8").Value = [double]"1.607e-5"

# F/G columns re-use the mean/stand.dev. strings already in the workbook
# (sum of effective temperatures "k" and developmental threshold "t")
$ws.Range("F4").Value = "929.354 ±49.111"
$ws.Range("F5").Value = "233.683 ±27.031"
$ws.Range("F6").Value = "243.945 ±45.301"
$ws.Range("F7").Value = "2602.996 ±297.464"
$ws.Range("F8").Value = "1207.431 ±489.288"

$ws.Range("G4").Value = "11.400 ±0.368"
$ws.Range("G5").Value = "15.437 ±0.305"
$ws.Range("G6").Value = "15.689 ±0.410"
$ws.Range("G7").Value = "9.375 ±0.846"
$ws.Range("G8").Value = "12.535 ±1.624"

$ws.Range("A2").Select()

# ------------------------------------------------------------------
# 2. dev.char sheet tidy-up
# ------------------------------------------------------------------
$dev = $wb.Worksheets.Item("dev.char")
$dev.Range("H8").Value = ""
$dev.Range("D14").Select()

# ------------------------------------------------------------------
# 3. mean.dev.time sheet selection
# ------------------------------------------------------------------
$mean = $wb.Worksheets.Item("mean.dev.time")
$mean.Range("E24").Select()

# ------------------------------------------------------------------
# 4. Workbook view: scroll tab strip so "dev.char" is first visible
# ------------------------------------------------------------------
$wb.Windows.Item(1).ScrollWorkbookTabs(1)
